$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-02 Wednesday" "2025-07-03 Thursday"

Replace-Text "167×6=1002" "809×9=7281"
Replace-Text "155×8=1240" "778×4=3112"
Replace-Text "549×9=4941" "949×4=3796"
Replace-Text "978×3=2934" "111×8=888"
Replace-Text "783×9=7047" "109×9=981"
Replace-Text "374×8=2992" "877×5=4385"
Replace-Text "703×2=1406" "114×7=798"
Replace-Text "556×7=3892" "315×9=2835"
Replace-Text "648×4=2592" "982×9=8838"
Replace-Text "638×9=5742" "552×6=3312"
Replace-Text "207×6=1242" "786×2=1572"
Replace-Text "978×7=6846" "396×8=3168"
Replace-Text "143×2=286" "442×9=3978"
Replace-Text "392×8=3136" "456×4=1824"
Replace-Text "605×9=5445" "365×4=1460"
Replace-Text "749×8=5992" "855×3=2565"
Replace-Text "548×4=2192" "992×2=1984"
Replace-Text "731×3=2193" "619×7=4333"
Replace-Text "782×9=7038" "955×2=1910"
Replace-Text "607×8=4856" "819×6=4914"
Replace-Text "682×7=4774" "912×7=6384"
Replace-Text "635×9=5715" "842×3=2526"
Replace-Text "378×5=1890" "925×4=3700"
Replace-Text "298×8=2384" "313×6=1878"
Replace-Text "720×4=2880" "908×8=7264"
